$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.730.21"
$ws.Range("E2").Value = "  +6.21%  "
$ws.Range("D3").Value = "2.046.52"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'251.44"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("D6").Value = "'0.650"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +17.08%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'59.59"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'0.375"
$ws.Range("E10").Value = "  +4.83%  "
$ws.Range("E11").Value = "  +3.88%  "
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").Value = "'0.903"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").Value = "'15.09"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("D15").Value = "2.342.43"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "'5.57"
$ws.Range("E16").Value = "  +6.36%  "
$ws.Range("D17").Value = "'20.51"
$ws.Range("E17").Value = "  +19.82%  "
$ws.Range("D18").Value = "2.044.17"
$ws.Range("E18").Value = "  +3.65%  "
$ws.Range("D19").Value = "37.591.46"
$ws.Range("E19").Value = "  +6.36%  "
$ws.Range("D20").Value = "'73.30"
$ws.Range("E20").Value = "  +4.78%  "
$ws.Range("D21").Value = "0.0₃0873"
$ws.Range("E21").Value = "  +4.68%  "
$ws.Range("D22").Value = "'5.34"
$ws.Range("E22").Value = "  +6.40%  "
$ws.Range("D23").Value = "'237.56"
$ws.Range("E23").Value = "  +2.45%  "
$ws.Range("D24").Value = "'2.67"
$ws.Range("E24").Value = "  +18.52%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").Value = "'9.57"
$ws.Range("E27").Value = "  +5.74%  "
$ws.Range("D28").Value = "'164.69"
$ws.Range("E28").Value = "  +1.02%  "
$ws.Range("D29").Value = "'19.92"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").Value = "'0.122"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("D31").Value = "'5.21"
$ws.Range("E31").Value = "  +9.26%  "
$ws.Range("D32").Value = "'0.113"
$ws.Range("E32").Value = "  +27.08%  "
$ws.Range("D33").Value = "'1.20"
$ws.Range("E33").Value = "  +7.33%  "
$ws.Range("D34").Value = "'4.72"
$ws.Range("E34").Value = "  +11.16%  "
$ws.Range("D35").Value = "'0.0613"
$ws.Range("E35").Value = "  +5.31%  "
$ws.Range("E36").Value = "  +7.84%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'6.11"
$ws.Range("E38").Value = "  +25.89%  "
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.83"
$ws.Range("E39").Value = "  +2.50%  "
$ws.Range("E40").Value = "  +16.02%  "
$ws.Range("D41").Value = "'1.23"
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("D42").Value = "'2.77"
$ws.Range("E42").Value = "  +23.53%  "
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").Value = "'0.0219"
$ws.Range("E44").Value = "  +5.23%  "
$ws.Range("E45").Value = "  +5.97%  "
$ws.Range("D46").Value = "'8.07"
$ws.Range("E46").Value = "  +8.89%  "
$ws.Range("D47").Value = "'16.88"
$ws.Range("E47").Value = "  +10.10%  "
$ws.Range("D48").Value = "'95.07"
$ws.Range("E48").Value = "  +5.29%  "
$ws.Range("D49").Value = "1.424.67"
$ws.Range("E49").Value = "  +3.92%  "
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "'47.46"
$ws.Range("E51").Value = "  +4.53%  "
